# Update cryptocurrency price/volume data per the Tue Feb 28 18:24:21 UTC 2023 GitHub Actions refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.461.41"
$ws.Range("E2").Value = "  +0.56%  "
$ws.Range("D3").Value = "1.640.83"
$ws.Range("E3").Value = "  +0.64%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.36%  "
$ws.Range("D5").Value = "'1.000"
$ws.Range("E5").Value = "  -0.40%  "
$ws.Range("D6").Value = "'303.56"
$ws.Range("E6").Value = "  +0.25%  "
$ws.Range("D7").Value = "'0.3796"
$ws.Range("E7").Value = "  +0.76%  "
$ws.Range("D8").Value = "'52.16"
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "'0.3615"
$ws.Range("E9").Value = "  +0.13%  "
$ws.Range("D10").Value = "'0.08170"
$ws.Range("E10").Value = "  +1.28%  "
$ws.Range("D11").Value = "'1.233"
$ws.Range("E11").Value = "  +0.45%  "
$ws.Range("D12").Value = "'1.003"
$ws.Range("E12").Value = "  -0.31%  "
$ws.Range("D13").Value = "'22.51"
$ws.Range("E13").Value = "  -0.48%  "
$ws.Range("D14").Value = "'6.458"
$ws.Range("E14").Value = "  -1.91%  "
$ws.Range("D15").Value = "'7.355"
$ws.Range("E15").Value = "  +1.69%  "
$ws.Range("D16").Value = "'0.00001238"
$ws.Range("E16").Value = "  -0.69%  "
$ws.Range("D17").Value = "1.637.67"
$ws.Range("E17").Value = "  +0.37%  "
$ws.Range("D18").Value = "'95.18"
$ws.Range("E18").Value = "  +1.66%  "
$ws.Range("D19").Value = "'0.06959"
$ws.Range("E19").Value = "  +0.64%  "
$ws.Range("D20").Value = "'17.55"
$ws.Range("E20").Value = "  -2.29%  "
$ws.Range("D21").Value = "'6.569"
$ws.Range("E21").Value = "  +1.38%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "  -0.40%  "
$ws.Range("D23").Value = "'12.51"
$ws.Range("E23").Value = "  -2.04%  "
$ws.Range("D24").Value = "23.456.64"
$ws.Range("E24").Value = "  +0.48%  "
$ws.Range("D25").Value = "'2.519"
$ws.Range("E25").Value = "  +3.50%  "
$ws.Range("D26").Value = "'3.056"
$ws.Range("E26").Value = "  -4.43%  "
$ws.Range("D27").Value = "'21.19"
$ws.Range("E27").Value = "  +0.57%  "
$ws.Range("D28").Value = "'151.92"
$ws.Range("D29").Value = "'5.256"
$ws.Range("E29").Value = "  -0.75%  "
$ws.Range("D30").Value = "'133.34"
$ws.Range("E30").Value = "  -1.05%  "
$ws.Range("D31").Value = "1.819.37"
$ws.Range("E31").Value = "  +0.17%  "
$ws.Range("D32").Value = "'1.096"
$ws.Range("E32").Value = "  +15.21%  "
$ws.Range("B33").Value = "WEMIXTOKEN"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D33").Value = "'2.156"
$ws.Range("E33").Value = "  -6.87%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "'6.583"
$ws.Range("E34").Value = "  -2.88%  "
$ws.Range("D35").Value = "'11.47"
$ws.Range("E35").Value = "  +5.30%  "
$ws.Range("D36").Value = "'0.02764"
$ws.Range("E36").Value = "  -2.59%  "
$ws.Range("D37").Value = "'0.2510"
$ws.Range("E37").Value = "  -0.99%  "
$ws.Range("D38").Value = "'0.08749"
$ws.Range("E38").Value = "  -1.03%  "
$ws.Range("D39").Value = "'5.989"
$ws.Range("E39").Value = "  -2.87%  "
$ws.Range("D40").Value = "'0.07044"
$ws.Range("E40").Value = "  -1.85%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'1.352"
$ws.Range("E41").Value = "  -0.91%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "'0.7045"
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("D43").Value = "'12.27"
$ws.Range("E43").Value = "  -0.57%  "
$ws.Range("D44").Value = "'15.67"
$ws.Range("E44").Value = "  -3.71%  "
$ws.Range("D45").Value = "'0.6537"
$ws.Range("E45").Value = "  +0.89%  "
$ws.Range("D46").Value = "'0.9994"
$ws.Range("E46").Value = "  -0.24%  "
$ws.Range("D47").Value = "'2.290"
$ws.Range("E47").Value = "  -1.52%  "
$ws.Range("D48").Value = "'3.964"
$ws.Range("E48").Value = "  -0.38%  "
$ws.Range("D49").Value = "'0.07984"
$ws.Range("E49").Value = "  +0.18%  "
$ws.Range("D50").Value = "'129.02"
$ws.Range("E50").Value = "  +1.79%  "
$ws.Range("D51").Value = "'1.193"
$ws.Range("E51").Value = "  -1.30%  "
